# 4 mdelo melhores rstds
# Reorders the model rows (column A) and refreshes the metric columns (B:I)
# for rows 2..26 of Sheet1 to reflect the 4 new "melhores" models added to
# the results table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New model name for every row (2..26), in order.
$names = @(
    "model_11_6_0",
    "model_11_6_22",
    "model_11_6_21",
    "model_11_6_20",
    "model_11_6_19",
    "model_11_6_18",
    "model_11_6_17",
    "model_11_6_16",
    "model_11_6_15",
    "model_11_6_14",
    "model_11_6_13",
    "model_11_6_23",
    "model_11_6_12",
    "model_11_6_10",
    "model_11_6_9",
    "model_11_6_8",
    "model_11_6_7",
    "model_11_6_6",
    "model_11_6_5",
    "model_11_6_4",
    "model_11_6_3",
    "model_11_6_2",
    "model_11_6_1",
    "model_11_6_11",
    "model_11_6_24"
)

# New metric values (r2, r2_test, r2_val, r2_vt, mse, mse_test, mse_val, mse_vt)
# -- identical for every row after the edit.
$metrics = @(
    0.3494677884409869,
    0.4069518043360018,
    0.1142347282625567,
    0.3450633061518301,
    0.7199474573135376,
    0.9062172770500183,
    0.8034555315971375,
    0.8578587174415588
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $names[$i]
    $ws.Range("B$row").Value = $metrics[0]
    $ws.Range("C$row").Value = $metrics[1]
    $ws.Range("D$row").Value = $metrics[2]
    $ws.Range("E$row").Value = $metrics[3]
    $ws.Range("F$row").Value = $metrics[4]
    $ws.Range("G$row").Value = $metrics[5]
    $ws.Range("H$row").Value = $metrics[6]
    $ws.Range("I$row").Value = $metrics[7]
}
